# Auto-generated Excel COM-interop script
# Applies cached market-price data updates across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as produced by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 170.42105
$ws.Range("I12").Value = 155.57143
$ws.Range("J12").Value = 212
$ws.Range("K12").Value = 155.57143
$ws.Range("L12").Value = 212
$ws.Range("M12").Value = 14.42857000000001
$ws.Range("N12").Value = -552
# Row 33
$ws.Range("H33").Value = 422.05884
$ws.Range("I33").Value = 385.9375
$ws.Range("K33").Value = 385.9375
$ws.Range("M33").Value = -156.9375
# Row 69
$ws.Range("H69").Value = 8512.857
$ws.Range("J69").Value = 8936.923000000001
$ws.Range("L69").Value = 26810.769
$ws.Range("N69").Value = -28558.769
# Row 72
$ws.Range("H72").Value = 8512.857
$ws.Range("J72").Value = 8936.923000000001
$ws.Range("L72").Value = 80432.307
$ws.Range("N72").Value = -89168.307
# Row 125
$ws.Range("H125").Value = 1122.1
$ws.Range("J125").Value = 1032.1111
$ws.Range("L125").Value = 9288.999900000001
$ws.Range("N125").Value = -14208.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 11731.444
$ws.Range("I2").Value = 697.875
$ws.Range("K2").Value = 697.875
$ws.Range("M2").Value = -584.875
# Row 61
$ws.Range("H61").Value = 5764.9165
$ws.Range("I61").Value = 5764.9165
$ws.Range("K61").Value = 5764.9165
$ws.Range("M61").Value = -5552.9165
# Row 93
$ws.Range("H93").Value = 34500
$ws.Range("J93").Value = 34500
$ws.Range("L93").Value = 34500
$ws.Range("N93").Value = -39492
# Row 116
$ws.Range("H116").Value = 11731.444
$ws.Range("I116").Value = 697.875
$ws.Range("K116").Value = 697.875
$ws.Range("M116").Value = 1596.125
# Row 136
$ws.Range("H136").Value = 5764.9165
$ws.Range("I136").Value = 5764.9165
$ws.Range("K136").Value = 17294.7495
$ws.Range("M136").Value = -14744.7495

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 11731.444
$ws.Range("I3").Value = 697.875
$ws.Range("K3").Value = 697.875
$ws.Range("M3").Value = -583.875
# Row 92
$ws.Range("H92").Value = 37500
$ws.Range("J92").Value = 37500
$ws.Range("L92").Value = 37500
$ws.Range("N92").Value = -42492
# Row 95
$ws.Range("H95").Value = 14250
$ws.Range("J95").Value = 14250
$ws.Range("L95").Value = 14250
$ws.Range("N95").Value = -19742
# Row 96
$ws.Range("H96").Value = 20729.5
$ws.Range("I96").Value = 15094.25
$ws.Range("J96").Value = 32000
$ws.Range("K96").Value = 15094.25
$ws.Range("L96").Value = 32000
$ws.Range("M96").Value = -12348.25
$ws.Range("N96").Value = -37492
# Row 106
$ws.Range("H106").Value = 68806.25
$ws.Range("J106").Value = 68806.25
$ws.Range("L106").Value = 68806.25
$ws.Range("N106").Value = -71330.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 37374.934
$ws.Range("I31").Value = 1352.0834
$ws.Range("K31").Value = 1352.0834
$ws.Range("M31").Value = -1057.0834
# Row 34
$ws.Range("H34").Value = 37374.934
$ws.Range("I34").Value = 1352.0834
$ws.Range("K34").Value = 1352.0834
$ws.Range("M34").Value = -1150.0834
# Row 35
$ws.Range("H35").Value = 5941.4375
$ws.Range("I35").Value = 6218.7856
$ws.Range("J35").Value = 4000
$ws.Range("K35").Value = 6218.7856
$ws.Range("L35").Value = 4000
$ws.Range("M35").Value = -5924.7856
$ws.Range("N35").Value = -4588
# Row 59
$ws.Range("H59").Value = 36411.2
$ws.Range("I59").Value = 5799
$ws.Range("J59").Value = 44064.25
$ws.Range("K59").Value = 5799
$ws.Range("L59").Value = 44064.25
$ws.Range("M59").Value = -4654
$ws.Range("N59").Value = -46354.25
# Row 95
$ws.Range("H95").Value = 44382.75
$ws.Range("J95").Value = 44382.75
$ws.Range("L95").Value = 44382.75
$ws.Range("N95").Value = -49874.75
# Row 134
$ws.Range("H134").Value = 6220.7144
$ws.Range("I134").Value = 3040.6316
$ws.Range("K134").Value = 9121.8948
$ws.Range("M134").Value = -6586.8948
# Row 137
$ws.Range("H137").Value = 78514
$ws.Range("J137").Value = 78514
$ws.Range("L137").Value = 78514
$ws.Range("N137").Value = -88714

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 141
$ws.Range("H141").Value = 5720.25
$ws.Range("I141").Value = 2887.875
$ws.Range("K141").Value = 8663.625
$ws.Range("M141").Value = -3483.625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1070.4375
$ws.Range("I113").Value = 915.2
$ws.Range("J113").Value = 1329.1666
$ws.Range("K113").Value = 915.2
$ws.Range("L113").Value = 1329.1666
$ws.Range("M113").Value = 1254.8
$ws.Range("N113").Value = -5669.1666
# Row 132
$ws.Range("H132").Value = 225021.97
$ws.Range("I132").Value = 245208.98
$ws.Range("K132").Value = 735626.9400000001
$ws.Range("M132").Value = -733096.9400000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7017.516
$ws.Range("I7").Value = 4389.4585
$ws.Range("J7").Value = 16028
$ws.Range("K7").Value = 4389.4585
$ws.Range("L7").Value = 16028
$ws.Range("M7").Value = -4277.4585
$ws.Range("N7").Value = -16252
# Row 46
$ws.Range("H46").Value = 2481.353
$ws.Range("I46").Value = 799.8
$ws.Range("J46").Value = 3182
$ws.Range("K46").Value = 799.8
$ws.Range("L46").Value = 3182
$ws.Range("M46").Value = -611.8
$ws.Range("N46").Value = -3558
# Row 61
$ws.Range("H61").Value = 5530.7646
$ws.Range("I61").Value = 3470.6
$ws.Range("K61").Value = 3470.6
$ws.Range("M61").Value = -3268.6
# Row 82
$ws.Range("H82").Value = 9646.294
$ws.Range("J82").Value = 9841.416999999999
$ws.Range("L82").Value = 9841.416999999999
$ws.Range("N82").Value = -10563.417
# Row 85
$ws.Range("H85").Value = 9646.294
$ws.Range("J85").Value = 9841.416999999999
$ws.Range("L85").Value = 9841.416999999999
$ws.Range("N85").Value = -12337.417
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
# Row 100
$ws.Range("H100").Value = 13424.167
$ws.Range("I100").Value = 10511
$ws.Range("K100").Value = 10511
$ws.Range("M100").Value = -9970
# Row 113
$ws.Range("H113").Value = 5530.7646
$ws.Range("I113").Value = 3470.6
$ws.Range("K113").Value = 3470.6
$ws.Range("M113").Value = -1300.6
# Row 126
$ws.Range("H126").Value = 7017.516
$ws.Range("I126").Value = 4389.4585
$ws.Range("J126").Value = 16028
$ws.Range("K126").Value = 13168.3755
$ws.Range("L126").Value = 48084
$ws.Range("M126").Value = -10698.3755
$ws.Range("N126").Value = -53024

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 627.06665
$ws.Range("I113").Value = 289.8889
$ws.Range("K113").Value = 869.6667
$ws.Range("M113").Value = 1300.3333
